$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.076.20'
$ws.Range("E2").Value = '  +0.01%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.368.80'
$ws.Range("E3").Value = '  +1.36%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("E5").Value = '  +0.21%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '95.18'
$ws.Range("E6").Value = '  +0.50%  '
$ws.Range("B7").Value = 'USDC'
$ws.Range("C7").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("B8").Value = 'XRP'
$ws.Range("C8").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.502'
$ws.Range("E8").Value = '  -0.46%  '
$ws.Range("E9").Value = '  -3.20%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.34'
$ws.Range("E10").Value = '  +0.37%  '
$ws.Range("E11").Value = '  +0.09%  '
$ws.Range("E12").Value = '  +1.01%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.59'
$ws.Range("E13").Value = '  -2.21%  '
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.736.62'
$ws.Range("E14").Value = '  +1.34%  '
$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.72'
$ws.Range("E15").Value = '  +0.08%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.366.53'
$ws.Range("E16").Value = '  +0.91%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.102.30'
$ws.Range("E18").Value = '  +0.22%  '
$ws.Range("E19").Value = '  -1.01%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.28'
$ws.Range("E20").Value = '  +1.80%  '
$ws.Range("E21").Value = '  -0.66%  '
$ws.Range("E22").Value = '  +0.07%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '235.48'
$ws.Range("E23").Value = '  -0.46%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.19'
$ws.Range("E24").Value = '  -2.61%  '
$ws.Range("E26").Value = '  -0.17%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.53'
$ws.Range("E27").Value = '  -0.63%  '
$ws.Range("E28").Value = '  +15.33%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.34'
$ws.Range("E29").Value = '  +2.28%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '32.35'
$ws.Range("E30").Value = '  +2.09%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.999'
$ws.Range("E31").Value = '  -0.11%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.02'
$ws.Range("E32").Value = '  +0.12%  '
$ws.Range("E33").Value = '  -1.07%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0728'
$ws.Range("E34").Value = '  +3.67%  '
$ws.Range("E35").Value = '  +6.31%  '
$ws.Range("E36").Value = '  +0.76%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '126.67'
$ws.Range("E37").Value = '  -9.58%  '
$ws.Range("E38").Value = '  -1.26%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.84'
$ws.Range("E39").Value = '  +2.90%  '
$ws.Range("E40").Value = '  -1.27%  '
$ws.Range("E41").Value = '  -1.07%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '20.81'
$ws.Range("E42").Value = '  -6.96%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.935.16'
$ws.Range("E43").Value = '  -0.40%  '
$ws.Range("E44").Value = '  -0.20%  '
$ws.Range("E45").Value = '  +4.09%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.32'
$ws.Range("E46").Value = '  -8.45%  '
$ws.Range("E47").Value = '  -0.67%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.596.31'
$ws.Range("E48").Value = '  +1.26%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.51'
$ws.Range("E49").Value = '  +2.35%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '71.81'
$ws.Range("E50").Value = '  -1.12%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.13'
$ws.Range("E51").Value = '  +0.74%  '
